# Werkbestand Projectadministratie.xlsx
# Commit: created an overview in another sheet within werkbestand for me to
# easily work on closing project.
#
# Concretely: the second sheet ("Projecten Afsluiten") is turned into a
# plain "Sheet2", its dynamic-array overview formula (row 3) is removed,
# and the first column header becomes "Projectnummer" instead of "Project".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Projecten Afsluiten")

# Rename the second sheet back to the generic default name.
$ws2.Name = "Sheet2"

# The header in A2 becomes "Projectnummer" (a brand-new shared string)
# instead of reusing "Project".
$ws2.Range("A2").Value = "Projectnummer"

# Drop the LET/FILTER array-formula row that built the "closing projects"
# overview - row 3 on Sheet2 - shrinking the used range back to A1:E2.
$ws2.Rows.Item(3).Delete()

# Restore cursor positions on each sheet. Select Sheet2's cell first so
# that selecting Sheet1's cell afterwards leaves Sheet1 as the active tab.
$ws2.Range("E18").Select() | Out-Null
$ws1.Range("C34").Select() | Out-Null
